$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NormalStage")

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 6

$ws.Range("B8").Select()
